# Apply updated cryptocurrency price/volume data as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.525.87"
$ws.Range("E2").Value = "'  +3.01%  "
$ws.Range("D3").Value = "'2.320.14"
$ws.Range("E3").Value = "'  +2.27%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("D5").Value = "'311.97"
$ws.Range("E5").Value = "'  +1.57%  "
$ws.Range("D6").Value = "'102.63"
$ws.Range("E6").Value = "'  +5.47%  "
$ws.Range("D7").Value = "'0.537"
$ws.Range("E7").Value = "'  +2.09%  "
$ws.Range("E8").Value = "'  -0.02%  "
$ws.Range("E9").Value = "'  +7.62%  "
$ws.Range("D10").Value = "'36.12"
$ws.Range("E10").Value = "'  +2.13%  "
$ws.Range("E11").Value = "'  +3.37%  "
$ws.Range("E12").Value = "'  -0.22%  "
$ws.Range("D14").Value = "'2.679.70"
$ws.Range("E14").Value = "'  +2.28%  "
$ws.Range("D15").Value = "'15.02"
$ws.Range("E15").Value = "'  +1.68%  "
$ws.Range("D16").Value = "'2.312.16"
$ws.Range("E16").Value = "'  +3.00%  "
$ws.Range("D17").Value = "'0.814"
$ws.Range("E17").Value = "'  +2.27%  "
$ws.Range("D18").Value = "'43.425.07"
$ws.Range("E18").Value = "'  +3.10%  "
$ws.Range("D19").Value = "'12.51"
$ws.Range("E19").Value = "'  +0.46%  "
$ws.Range("D20").Value = "'0.0₃0928"
$ws.Range("E20").Value = "'  +2.25%  "
$ws.Range("D21").Value = "'6.16"
$ws.Range("E21").Value = "'  +2.17%  "
$ws.Range("D22").Value = "'68.43"
$ws.Range("E22").Value = "'  +0.17%  "
$ws.Range("D23").Value = "'242.31"
$ws.Range("E23").Value = "'  +1.63%  "
$ws.Range("E24").Value = "'  +5.50%  "
$ws.Range("E25").Value = "'  +2.49%  "
$ws.Range("D26").Value = "'0.997"
$ws.Range("E26").Value = "'  -0.28%  "
$ws.Range("E27").Value = "'  -1.22%  "
$ws.Range("D28").Value = "'24.78"
$ws.Range("E28").Value = "'  +4.67%  "
$ws.Range("D29").Value = "'37.69"
$ws.Range("E29").Value = "'  -0.09%  "
$ws.Range("D30").Value = "'9.66"
$ws.Range("E30").Value = "'  +1.64%  "
$ws.Range("D32").Value = "'167.22"
$ws.Range("E32").Value = "'  +3.22%  "
$ws.Range("D33").Value = "'5.34"
$ws.Range("E33").Value = "'  +1.86%  "
$ws.Range("E34").Value = "'  +0.06%  "
$ws.Range("D35").Value = "'2.52"
$ws.Range("E35").Value = "'  +6.51%  "
$ws.Range("D36").Value = "'0.0748"
$ws.Range("E36").Value = "'  +1.29%  "
$ws.Range("E37").Value = "'  -2.38%  "
$ws.Range("E38").Value = "'  +3.25%  "
$ws.Range("E40").Value = "'  +2.44%  "
$ws.Range("E41").Value = "'  +1.88%  "
$ws.Range("E42").Value = "'  +7.83%  "
$ws.Range("D43").Value = "'19.80"
$ws.Range("E43").Value = "'  +4.63%  "
$ws.Range("D44").Value = "'2.32"
$ws.Range("E44").Value = "'  +0.20%  "
$ws.Range("B45").Value = "'VeChain"
$ws.Range("C45").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0291"
$ws.Range("E45").Value = "'  +2.96%  "
$ws.Range("B46").Value = "'Maker"
$ws.Range("C46").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "'1.982.23"
$ws.Range("E46").Value = "'  +2.04%  "
$ws.Range("E47").Value = "'  +4.33%  "
$ws.Range("D48").Value = "'9.87"
$ws.Range("E48").Value = "'  -1.24%  "
$ws.Range("D49").Value = "'56.03"
$ws.Range("E49").Value = "'  +4.66%  "
$ws.Range("D50").Value = "'2.95"
$ws.Range("E50").Value = "'  +12.20%  "
$ws.Range("E51").Value = "'  +6.51%  "
